# Applies the change: insert a new weekly price record as row 87
# (pushing the former rows 87-191 down to 88-192), and populate the
# new row 87 with the same record data as the row that used to be at
# position 87, except for an updated Fecha (D) and Volumen (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 87; Excel shifts rows 87:191 down to 88:192
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C87").Value = 'Los Lagos'
$ws.Range("D87").Value = 44467
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = 100112023
$ws.Range("G87").Value = 'Brócoli'
$ws.Range("H87").Value = 'Sin especificar'
$ws.Range("I87").Value = 'Primera'
$ws.Range("J87").Value = 1400
$ws.Range("K87").Value = 1300
$ws.Range("L87").Value = 1300
$ws.Range("M87").Value = 1300
$ws.Range("N87").Value = '$/unidad'
$ws.Range("O87").Value = 'Región Metropolitana'
$ws.Range("P87").Value = 1300
$ws.Range("Q87").Value = 1
$ws.Range("R87").Value = 'Hortaliza'

# Make sure the date cell keeps a date number format like the rest of column D
$ws.Range("D87").NumberFormat = $ws.Range("D88").NumberFormat
